# "fixed path in excel files"
#
# 1) Sheet "protocoltestcasedetails": two test cases (testcase25_oracle_oracle_etljob
#    and testcase26_names_fullname_etljob) were removed entirely (whole rows 4 and 5
#    deleted, shifting the remaining cases up); the CONCAT() formulas that build the
#    test-case file path were switched from backslashes to forward slashes.
# 2) Sheet "protocol": the results-path value was switched from backslashes to
#    forward slashes as well.
# 3) Selection (active cell) bookkeeping left behind by the editing session.

$wb = $excel.ActiveWorkbook

$wsProtocol = $wb.Worksheets.Item("protocol")
$wsCases    = $wb.Worksheets.Item("protocoltestcasedetails")

# --- protocoltestcasedetails: drop the two retired test cases -------------
$wsCases.Rows("4:5").Delete()

# --- fix the slash direction in the CONCAT formulas that remain -----------
$wsCases.Range("C2").Formula = '=_xlfn.CONCAT("test/testcases/",B2,".xlsx")'
$wsCases.Range("C3").Formula = '=_xlfn.CONCAT("test/testcases/",B3,".xlsx")'
$wsCases.Range("C4").Formula = '=_xlfn.CONCAT("test/testcases/",B4,".xlsx")'
$wsCases.Range("C5").Formula = '=_xlfn.CONCAT("test/testcases/",B5,".xlsx")'
$wsCases.Range("C6").Formula = '=_xlfn.CONCAT("test/testcases/",B6,".xlsx")'
$wsCases.Range("C7").Formula = '=_xlfn.CONCAT("test/testcases/",B7,".xlsx")'

# --- protocol: fix the slash direction in the results path ----------------
$wsProtocol.Range("B3").Value = "test/results/"

# --- restore the selections left by the editing session -------------------
$wsProtocol.Range("B9").Select()
$wsCases.Range("B11").Select()
